$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "temp solve of RWheel" - set Fitness column (C2:C12) to a fixed temp value
$ws.Range("C2:C12").Value = 3974
